$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1051
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 36
$ws.Range("H2").Value = 19
$ws.Range("I2").Value = 19
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1770
$ws.Range("L2").Value = 598
$ws.Range("M2").Value = 1171
$ws.Range("N2").Value = 1003
$ws.Range("O2").Value = 169
$ws.Range("P2").Value = 87
$ws.Range("Q2").Value = -93
$ws.Range("R2").Value = -108
$ws.Range("S2").Value = 198
$ws.Range("T2").Value = 130
$ws.Range("U2").Value = -223
$ws.Range("V2").Value = 342
$ws.Range("W2").Value = 1.08
$ws.Range("X2").Value = 1.82
$ws.Range("Y2").Value = 1.86
$ws.Range("Z2").Value = 1.18
$ws.Range("AA2").Value = 51.08
$ws.Range("AB2").Value = 1056.65
$ws.Range("AC2").Value = 108
$ws.Range("AD2").Value = 82.40000000000001
$ws.Range("AE2").Value = 5775
$ws.Range("AF2").Value = 1.54
$ws.Range("AG2").Value = 170
$ws.Range("AH2").Value = 1.91
$ws.Range("AI2").Value = 157.76
$ws.Range("AJ2").Value = 17360000

# Row 3
$ws.Range("D3").Value = 1062
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 199
$ws.Range("H3").Value = 140
$ws.Range("I3").Value = 137
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1947
$ws.Range("L3").Value = 671
$ws.Range("M3").Value = 1277
$ws.Range("N3").Value = 1103
$ws.Range("O3").Value = 174
$ws.Range("P3").Value = 87
$ws.Range("Q3").Value = -15
$ws.Range("R3").Value = 271
$ws.Range("S3").Value = 62
$ws.Range("T3").Value = 30
$ws.Range("U3").Value = -45
$ws.Range("V3").Value = 414
$ws.Range("W3").Value = 1.48
$ws.Range("X3").Value = 13.16
$ws.Range("Y3").Value = 13
$ws.Range("Z3").Value = 7.52
$ws.Range("AA3").Value = 52.52
$ws.Range("AB3").Value = 1172.51
$ws.Range("AC3").Value = 789
$ws.Range("AD3").Value = 9.460000000000001
$ws.Range("AE3").Value = 6355
$ws.Range("AF3").Value = 1.17
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 2.68
$ws.Range("AI3").Value = 25.36
$ws.Range("AJ3").Value = 17360000

# Row 4
$ws.Range("D4").Value = 1084
$ws.Range("E4").Value = -3
$ws.Range("F4").Value = -3
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = -3
$ws.Range("I4").Value = -4
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1813
$ws.Range("L4").Value = 579
$ws.Range("M4").Value = 1234
$ws.Range("N4").Value = 1064
$ws.Range("O4").Value = 170
$ws.Range("P4").Value = 87
$ws.Range("Q4").Value = -145
$ws.Range("R4").Value = -24
$ws.Range("S4").Value = -63
$ws.Range("T4").Value = 22
$ws.Range("U4").Value = -167
$ws.Range("V4").Value = 345
$ws.Range("W4").Value = -0.26
$ws.Range("X4").Value = -0.23
$ws.Range("Y4").Value = -0.35
$ws.Range("Z4").Value = -0.13
$ws.Range("AA4").Value = 46.95
$ws.Range("AB4").Value = 1126.92
$ws.Range("AC4").Value = -22
$ws.Range("AD4").Value = -310.61
$ws.Range("AE4").Value = 6129
$ws.Range("AF4").Value = 1.11
$ws.Range("AG4").Value = 220
$ws.Range("AH4").Value = 3.24
$ws.Range("AI4").Value = -1004.9
$ws.Range("AJ4").Value = 17360000

# Row 5
$ws.Range("D5").Value = 983
$ws.Range("E5").Value = -28
$ws.Range("F5").Value = -28
$ws.Range("G5").Value = -197
$ws.Range("H5").Value = -156
$ws.Range("I5").Value = -157
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1899
$ws.Range("L5").Value = 531
$ws.Range("M5").Value = 1368
$ws.Range("N5").Value = 1158
$ws.Range("O5").Value = 210
$ws.Range("P5").Value = 87
$ws.Range("Q5").Value = 131
$ws.Range("R5").Value = -134
$ws.Range("S5").Value = 42
$ws.Range("T5").Value = 14
$ws.Range("U5").Value = 117
$ws.Range("V5").Value = 231
$ws.Range("W5").Value = -2.83
$ws.Range("X5").Value = -15.85
$ws.Range("Y5").Value = -14.17
$ws.Range("Z5").Value = -8.390000000000001
$ws.Range("AA5").Value = 38.83
$ws.Range("AB5").Value = 942.1900000000001
$ws.Range("AC5").Value = -907
$ws.Range("AD5").Value = -5.3
$ws.Range("AE5").Value = 6860
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 17360000

# Row 6
$ws.Range("D6").Value = 939
$ws.Range("E6").Value = -11
$ws.Range("F6").Value = -11
$ws.Range("G6").Value = -96
$ws.Range("H6").Value = -72
$ws.Range("I6").Value = -73
$ws.Range("K6").Value = 2110
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = 1310
$ws.Range("N6").Value = 1110
$ws.Range("P6").Value = 88
$ws.Range("Q6").Value = 54
$ws.Range("R6").Value = -487
$ws.Range("S6").Value = 305
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 44
$ws.Range("V6").Value = 426
$ws.Range("W6").Value = -1.16
$ws.Range("X6").Value = -7.67
$ws.Range("Y6").Value = -6.44
$ws.Range("Z6").Value = -3.59
$ws.Range("AA6").Value = 61.06
$ws.Range("AB6").Value = 859.05
$ws.Range("AC6").Value = -420
$ws.Range("AD6").Value = -11.39
$ws.Range("AE6").Value = 6505
$ws.Range("AF6").Value = 0.73
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 17543127

# Rows 7-9: clear all estimate data except columns A-C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()